$d = $word.ActiveDocument

# The "Requisitos" bullet list (one run per "<code> - <name>  (Requisito)" line,
# 26 lines total, each terminated by a manual line break) is reordered in
# place: same 26 items, new sequence, no text content changes. Every item
# string is unique across the whole document, so Find/Replace against
# $d.Content (re-evaluated fresh -- it always reflects the live document)
# can target each line unambiguously without needing to hand-roll Range math.
#
# Find/Replace rewrites text at its current position; it cannot relocate text
# to a different slot by itself. So the reorder is done in two passes:
#   Phase 1 - stamp each position that must change (original order) with a
#             unique placeholder token, so no original requirement text is
#             left around to ambiguously collide with during phase 2.
#   Phase 2 - rewrite each placeholder with whatever item text belongs in
#             that same position under the new order.
# A position whose item does not move (same text before/after) is left
# untouched entirely.

# Phase 1: stamp every position that changes with a placeholder token.
$d.Content.Find.Execute("LOB1003 -  Cálculo I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT00@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1004 -  Cálculo II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT01@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1006 -  Cálculo IV  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT02@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT03@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1012 -  Estatística  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT04@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1018 -  Física I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT05@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1019 -  Física II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT06@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1021 -  Física IV  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT07@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1036 -  Geometria Analítica  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT08@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1037 -  Àlgebra Linear  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT09@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1038 -  Física Experimental I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT10@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1039 -  Física Experimental III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT11@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1041 -  Física Experimental II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT12@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1042 -  Física Experimental IV  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT13@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT14@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1052 -  Cálculo III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT15@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1053 -  Física III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT16@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT17@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT18@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3218 -  Introdução à Engenharia Física  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT19@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3236 -  Processos de Fabricação  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT20@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3241 -  Química de Materiais  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT21@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3261 -  Métodos Numéricos e Aplicações  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT23@@", 2) | Out-Null
$d.Content.Find.Execute("LOQ4095 -  Química Geral Experimental  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT24@@", 2) | Out-Null
$d.Content.Find.Execute("LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT25@@", 2) | Out-Null

# Phase 2: each placeholder becomes the item text for its new (target) position.
$d.Content.Find.Execute("@@SLOT00@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1053 -  Física III  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT01@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT02@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT03@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1036 -  Geometria Analítica  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT04@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1037 -  Àlgebra Linear  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT05@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1041 -  Física Experimental II  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT06@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1042 -  Física Experimental IV  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT07@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4095 -  Química Geral Experimental  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT08@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1039 -  Física Experimental III  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT09@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1018 -  Física I  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT10@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT11@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1004 -  Cálculo II  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT12@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1038 -  Física Experimental I  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT13@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1052 -  Cálculo III  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT14@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3236 -  Processos de Fabricação  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT15@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3261 -  Métodos Numéricos e Aplicações  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT16@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3218 -  Introdução à Engenharia Física  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT17@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1003 -  Cálculo I  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT18@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1006 -  Cálculo IV  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT19@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3241 -  Química de Materiais  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT20@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1021 -  Física IV  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT21@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT23@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT24@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1012 -  Estatística  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT25@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1019 -  Física II  (Requisito)", 2) | Out-Null

Write-Host "Requisitos list reordered."